$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 4; $r -le 23; $r++) {
    $ws.Cells.Item($r, 3).Value = 3462
}
